$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# NOV-2020 tracker sheet: correct row 25's comment, and append six more daily
# rows (26-31) covering 2020-11-25 through 2020-11-30.
# ---------------------------------------------------------------------------

# Prime rows 26 & 27 with the same formatting as row 25 (the last filled
# data row: No/Date/Application/Comments/% of completion/Status/blank).
$ws.Range("A25:G25").Copy()
$ws.Range("A26:G26").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A27:G27").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# No. / Date for the two new filled rows.
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = 44160
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = 44161

# Application column for rows 26 & 27.
$ws.Range("C26").Value = "Sonia and nMVAR"
$ws.Range("C27").Value = "Sonia and nMVAR"

# Comments column - row 27 first, then fix up row 25's comment, then row 26,
# matching the order the strings were authored in.
$ws.Range("D27").Value = "Sanity testing on B2C app, QMVAR site, GSS site and Hayaai site. Regression testing on Sonia Application(Soukastu) nMVAR_AI and nMVAR_Search"
$ws.Range("D25").Value = "Regression testing on nMVAR_QA, nMVAR_Clct and nMVAR_Mnt"
$ws.Range("D26").Value = "Sanity testing on B2C app, QMVAR site, GSS site and Hayaai site. Regression testing on Sonia Application(IVC and Hosyou), nMVAR_QA, nMVAR_Clct and nMVAR_Mnt"

# % of completion / Status for rows 26 & 27.
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = "Completed"
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = "Completed"

# Rows 28-31: blank placeholder days (2020-11-27 .. 2020-11-30). Base the
# No./Date/Application/Comments formatting on row 25 (s=16 for C & D) then
# restyle the %/Status/blank columns to the unfilled look used by earlier
# blank rows (s=1) before stamping in the serial No. and Date values.
$ws.Range("A25:G25").Copy()
$ws.Range("A28:G28").PasteSpecial(-4122)
$ws.Range("A29:G29").PasteSpecial(-4122)
$ws.Range("A30:G30").PasteSpecial(-4122)
$ws.Range("A31:G31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E8:G8").Copy()
$ws.Range("E28:G28").PasteSpecial(-4122)
$ws.Range("E29:G29").PasteSpecial(-4122)
$ws.Range("E30:G30").PasteSpecial(-4122)
$ws.Range("E31:G31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A28").Value = 27
$ws.Range("B28").Value = 44162
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = 44163
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = 44164
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = 44165

# Match the workbook's on-screen scroll position / selection after the edit.
$ws.Application.ActiveWindow.ScrollRow = 21
$ws.Range("B30").Select()
